# Handles float input without breaking stuff
#
# The quiz was re-graded: score summary (rows 10-12) updated, the marking
# scheme's "wrong answer" penalty is now stored as a real number instead of
# text, the 3rd answer-key block (columns G:H) is dropped entirely, and the
# 2nd block (columns D:E) is trimmed down to just the first few rows - with
# some of its "revealed answer" cells promoted into column A instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Row 10/11/12 label cells (A10,A11,A12) adopt the "mtitleStyle" ----
#     (same style already used by the row-9 header cells) while keeping
#     their existing text ("No.", "Marking", "Total").
$ws.Range("A9").Copy()
$ws.Range("A10:A12").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- 2. Score-summary numbers -------------------------------------------
$ws.Range("B10").Value2 = 16
$ws.Range("D10").Value2 = 12
$ws.Range("E10").Value2 = 28

$ws.Range("B11").Value2 = 4
$ws.Range("C11").Value2 = -1          # now numeric, used to be the text "-1"

$ws.Range("B12").Value2 = 64
$ws.Range("E12").Value2 = "64/112"    # used to read "Absent"

# --- 3. Drop the whole 3rd answer-key block (columns G:H) ---------------
$ws.Range("G15:H40").Clear()

# --- 4. Trim the 2nd answer-key block (columns D:E) down to rows 16-18 --
#     Rows 19-40 lose their D/E "Student Ans"/"Correct Ans" cells outright.
$ws.Range("D19:E40").Clear()

# --- 5. Promote certain "revealed answer" cells from normalStyle (blank) -
#     to correctStyle, filling in the value that belongs there. Source the
#     formatting from an existing correctStyle cell (B10) so the shared
#     cellXf is reused instead of minting a new one.
$ws.Range("B10").Copy()
$fmtTargets = $ws.Range("D16,D17,A18,D18,A20,A21,A25,A26,A30,A31,A33,A34,A35,A36,A38,A39")
$fmtTargets.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("D16").Value2 = "Option A"
$ws.Range("D17").Value2 = "Option C"
$ws.Range("A18").Value2 = "Option B"
$ws.Range("D18").Value2 = "Option D"
$ws.Range("A20").Value2 = "Option B"
$ws.Range("A21").Value2 = "Option C"
$ws.Range("A25").Value2 = "Option A"
$ws.Range("A26").Value2 = "Option C"
$ws.Range("A30").Value2 = "Option B"
$ws.Range("A31").Value2 = "Option D"
$ws.Range("A33").Value2 = "Option D"
$ws.Range("A34").Value2 = "Option B"
$ws.Range("A35").Value2 = "Option D"
$ws.Range("A36").Value2 = "Option A"
$ws.Range("A38").Value2 = "Option A"
$ws.Range("A39").Value2 = "Option D"
